$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cycle = @(10002, 10003, 10004, 10005, 10006, 10007, 10008, 10009, 10010)

$startRow = 102
$endRow = 146
$deviceId = 3000121

for ($row = $startRow; $row -le $endRow; $row++) {
    $idx = ($row - $startRow) % 9
    $regCntrId = $cycle[$idx]

    $ws.Cells.Item($row, 1).Value = $regCntrId
    $ws.Cells.Item($row, 2).Value = $deviceId
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"

    $deviceId++
}

# Update the view: scroll so row 128 is the top-left visible row, and select A102:F146
$ws.Application.ActiveWindow.ScrollRow = 128
$ws.Range("A102:F146").Select()
